$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9888435254646026
$ws.Range("D2").Value = 0.1960285505503165
$ws.Range("E2").Value = 0.254058699034891
$ws.Range("F2").Value = 2.041934024148162
$ws.Range("G2").Value = 0.002510227643564386
$ws.Range("J2").Value = 0.3680183331625582
$ws.Range("K2").Value = 0.3337214568804541
$ws.Range("L2").Value = 0.1148595582085647
$ws.Range("M2").Value = 0.1991049508157339
$ws.Range("O2").Value = 5.111854132152075
$ws.Range("B3").Value = 0.9691841713211886
$ws.Range("D3").Value = 0.1961467046099656
$ws.Range("E3").Value = 0.2555050508993162
$ws.Range("F3").Value = 2.051232215889364
$ws.Range("G3").Value = 0.002512851120024611
$ws.Range("J3").Value = 0.3700752552627353
$ws.Range("K3").Value = 0.2915605766560532
$ws.Range("L3").Value = 0.1085986163030057
$ws.Range("M3").Value = 0.1933142848010725
$ws.Range("O3").Value = 5.137784986601332
$ws.Range("B4").Value = 0.9575015919575947
$ws.Range("D4").Value = 0.1962734576156038
$ws.Range("E4").Value = 0.2564495180953346
$ws.Range("F4").Value = 2.057812302379375
$ws.Range("G4").Value = 0.002514549278901908
$ws.Range("J4").Value = 0.3714069037410188
$ws.Range("K4").Value = 0.2656168504571212
$ws.Range("L4").Value = 0.1047847123493426
$ws.Range("M4").Value = 0.1898312177501964
$ws.Range("O4").Value = 5.155995221828078
$ws.Range("B5").Value = 0.9528390455586475
$ws.Range("D5").Value = 0.1963387881446863
$ws.Range("E5").Value = 0.2568486065458897
$ws.Range("F5").Value = 2.060712966378418
$ws.Range("G5").Value = 0.002515263317904631
$ws.Range("J5").Value = 0.3719668673353853
$ws.Range("K5").Value = 0.2550308866107684
$ws.Range("L5").Value = 0.1032382661623146
$ws.Range("M5").Value = 0.1884301877890806
$ws.Range("O5").Value = 5.163991899706161
$ws.Range("B6").Value = 0.9520707804224173
$ws.Range("D6").Value = 0.1963504637127116
$ws.Range("E6").Value = 0.256915733970704
$ws.Range("F6").Value = 2.061207866435282
$ws.Range("G6").Value = 0.002515383215778778
$ws.Range("J6").Value = 0.3720608950878068
$ws.Range("K6").Value = 0.253272286333555
$ws.Range("L6").Value = 0.1029819516550603
$ws.Range("M6").Value = 0.1881986600298688
$ws.Range("O6").Value = 5.165354533513693
$ws.Range("B7").Value = 0.9574383128922932
$ws.Range("D7").Value = 0.1962742832419409
$ws.Range("E7").Value = 0.2564548427671394
$ws.Range("F7").Value = 2.057850533828962
$ws.Range("G7").Value = 0.002514558819274704
$ws.Range("J7").Value = 0.3714143854896212
$ws.Range("K7").Value = 0.2654741390687718
$ws.Range("L7").Value = 0.1047638248915419
$ws.Range("M7").Value = 0.1898122484980753
$ws.Range("O7").Value = 5.156100735698061
$ws.Range("B8").Value = 0.9819846993280237
$ws.Range("D8").Value = 0.1960580676029267
$ws.Range("E8").Value = 0.2545457118673724
$ws.Range("F8").Value = 2.044959400640728
$ws.Range("G8").Value = 0.002511114131211626
$ws.Range("J8").Value = 0.3687133234936604
$ws.Range("K8").Value = 0.3191965106682346
$ws.Range("L8").Value = 0.1126945590790172
$ws.Range("M8").Value = 0.1970933767317256
$ws.Range("O8").Value = 5.120320315984202
$ws.Range("B9").Value = 1.033180689164169
$ws.Range("D9").Value = 0.1960622099243778
$ws.Range("E9").Value = 0.2512481637570971
$ws.Range("F9").Value = 2.026582046808727
$ws.Range("G9").Value = 0.002505049103377309
$ws.Range("J9").Value = 0.3639600223844472
$ws.Range("K9").Value = 0.4240745291840256
$ws.Range("L9").Value = 0.1284829461688162
$ws.Range("M9").Value = 0.2119412074431715
$ws.Range("O9").Value = 5.068299715067866
$ws.Range("B10").Value = 1.07263723059981
$ws.Range("D10").Value = 0.1963237889539187
$ws.Range("E10").Value = 0.2490957325793355
$ws.Range("F10").Value = 2.017277789833258
$ws.Range("G10").Value = 0.002501009624164248
$ws.Range("J10").Value = 0.3607968768074938
$ws.Range("K10").Value = 0.5008204582658209
$ws.Range("L10").Value = 0.1402219052314138
$ws.Range("M10").Value = 0.2231915708743557
$ws.Range("O10").Value = 5.041125984329256
$ws.Range("B11").Value = 1.090982636389043
$ws.Range("D11").Value = 0.1964984154201801
$ws.Range("E11").Value = 0.2481748477935435
$ws.Range("F11").Value = 2.013954466768595
$ws.Range("G11").Value = 0.002499261518042695
$ws.Range("J11").Value = 0.3594288709765854
$ws.Range("K11").Value = 0.5356634244135989
$ws.Range("L11").Value = 0.1455915613760226
$ws.Range("M11").Value = 0.2283827478440799
$ws.Range("O11").Value = 5.031159396204885
$ws.Range("B12").Value = 1.097986086158812
$ws.Range("D12").Value = 0.1965724968063398
$ws.Range("E12").Value = 0.247834482998134
$ws.Range("F12").Value = 2.012826570044027
$ws.Range("G12").Value = 0.002498612355301897
$ws.Range("J12").Value = 0.3589210074769156
$ws.Range("K12").Value = 0.5488471139563842
$ws.Range("L12").Value = 0.1476290538646623
$ws.Range("M12").Value = 0.2303589370510579
$ws.Range("O12").Value = 5.027729395498795
$ws.Range("B13").Value = 1.096475266587476
$ws.Range("D13").Value = 0.196556189019546
$ws.Range("E13").Value = 0.247907415475165
$ws.Range("F13").Value = 2.013063678238808
$ws.Range("G13").Value = 0.002498751595416261
$ws.Range("J13").Value = 0.3590299330751834
$ws.Range("K13").Value = 0.5460082515321005
$ws.Range("L13").Value = 0.1471900621898357
$ws.Range("M13").Value = 0.2299328689149789
$ws.Range("O13").Value = 5.028452806247202
$ws.Range("B14").Value = 1.0915576864397
$ws.Range("D14").Value = 0.1965043509891089
$ws.Range("E14").Value = 0.2481466784840212
$ws.Range("F14").Value = 2.013859058126258
$ws.Range("G14").Value = 0.002499207854610553
$ws.Range("J14").Value = 0.3593868850680724
$ws.Range("K14").Value = 0.5367482703105964
$ws.Range("L14").Value = 0.1457591052664071
$ws.Range("M14").Value = 0.2285451224347312
$ws.Range("O14").Value = 5.030870312544522
$ws.Range("B15").Value = 1.088552860431435
$ws.Range("D15").Value = 0.1964736332051018
$ws.Range("E15").Value = 0.2482943212377409
$ws.Range("F15").Value = 2.014363250684411
$ws.Range("G15").Value = 0.002499488992855511
$ws.Range("J15").Value = 0.3596068521500815
$ws.Range("K15").Value = 0.5310748668470069
$ws.Range("L15").Value = 0.1448831356116642
$ws.Range("M15").Value = 0.2276964382456228
$ws.Range("O15").Value = 5.032395913380014
$ws.Range("B16").Value = 1.071446244158921
$ws.Range("D16").Value = 0.1963134922601171
$ws.Range("E16").Value = 0.2491570849543239
$ws.Range("F16").Value = 2.017513258096031
$ws.Range("G16").Value = 0.002501125662634992
$ws.Range("J16").Value = 0.3608877039713789
$ws.Range("K16").Value = 0.4985419421441293
$ws.Range("L16").Value = 0.1398715691149022
$ws.Range("M16").Value = 0.2228537787250104
$ws.Range("O16").Value = 5.041825493958669
$ws.Range("B17").Value = 1.061053022677385
$ws.Range("D17").Value = 0.1962294672450682
$ws.Range("E17").Value = 0.2497012686864961
$ws.Range("F17").Value = 2.019678438513793
$ws.Range("G17").Value = 0.00250215258029072
$ws.Range("J17").Value = 0.3616916104281112
$ws.Range("K17").Value = 0.4785658792991114
$ws.Range("L17").Value = 0.1368046111036847
$ws.Range("M17").Value = 0.219901644680057
$ws.Range("O17").Value = 5.048223443804233
$ws.Range("B18").Value = 1.055112463749964
$ws.Range("D18").Value = 0.1961863774368098
$ws.Range("E18").Value = 0.2500197545899008
$ws.Range("F18").Value = 2.021009384386559
$ws.Range("G18").Value = 0.002502751661420996
$ws.Range("J18").Value = 0.3621606736253575
$ws.Range("K18").Value = 0.4670697065991476
$ws.Range("L18").Value = 0.1350433665052293
$ws.Range("M18").Value = 0.2182105658875741
$ws.Range("O18").Value = 5.052128807608113
$ws.Range("B19").Value = 1.053107524159174
$ws.Range("D19").Value = 0.1961726891955777
$ws.Range("E19").Value = 0.250128531384787
$ws.Range("F19").Value = 2.021474725448385
$ws.Range("G19").Value = 0.002502955948924291
$ws.Range("J19").Value = 0.3623206380932023
$ws.Range("K19").Value = 0.4631762085141418
$ws.Range("L19").Value = 0.1344475221954298
$ws.Range("M19").Value = 0.21763918689485
$ws.Range("O19").Value = 5.053489822539177
$ws.Range("B20").Value = 1.062155537636954
$ws.Range("D20").Value = 0.1962378699506147
$ws.Range("E20").Value = 0.2496427718157452
$ws.Range("F20").Value = 2.019439094225092
$ws.Range("G20").Value = 0.002502042391598324
$ws.Range("J20").Value = 0.3616053423263699
$ws.Range("K20").Value = 0.4806930395472477
$ws.Range("L20").Value = 0.1371308063481962
$ws.Range("M20").Value = 0.2202151901921638
$ws.Range("O20").Value = 5.047519041534343
$ws.Range("B21").Value = 1.093000571574663
$ws.Range("D21").Value = 0.1965193615668568
$ws.Range("E21").Value = 0.2480761745995883
$ws.Range("F21").Value = 2.013621893419469
$ws.Range("G21").Value = 0.002499073493079191
$ws.Range("J21").Value = 0.3592817638115497
$ws.Range("K21").Value = 0.5394684436606099
$ws.Range("L21").Value = 0.1461793011993677
$ws.Range("M21").Value = 0.228952456003725
$ws.Range("O21").Value = 5.030150894756758
$ws.Range("B22").Value = 1.113488315477298
$ws.Range("D22").Value = 0.1967496742375872
$ws.Range("E22").Value = 0.2471009966986237
$ws.Range("F22").Value = 2.010581019033992
$ws.Range("G22").Value = 0.002497207769348082
$ws.Range("J22").Value = 0.3578224401493824
$ws.Range("K22").Value = 0.5778195707460725
$ws.Range("L22").Value = 0.1521169724205294
$ws.Range("M22").Value = 0.2347233363745929
$ws.Range("O22").Value = 5.020805542657627
$ws.Range("B23").Value = 1.102523727748775
$ws.Range("D23").Value = 0.196622526861475
$ws.Range("E23").Value = 0.2476170209480886
$ws.Range("F23").Value = 2.012134416203452
$ws.Range("G23").Value = 0.002498196732637164
$ws.Range("J23").Value = 0.3585958947910131
$ws.Range("K23").Value = 0.557356736971883
$ws.Range("L23").Value = 0.1489457757474639
$ws.Range("M23").Value = 0.2316378132607895
$ws.Range("O23").Value = 5.02560989094232
$ws.Range("B24").Value = 1.061656982674691
$ws.Range("D24").Value = 0.1962340548341004
$ws.Range("E24").Value = 0.2496692007205974
$ws.Range("F24").Value = 2.01954703339706
$ws.Range("G24").Value = 0.00250209218082331
$ws.Range("J24").Value = 0.3616443226812862
$ws.Range("K24").Value = 0.4797313866970683
$ws.Range("L24").Value = 0.1369833272690784
$ws.Range("M24").Value = 0.220073417122876
$ws.Range("O24").Value = 5.047836794381254
$ws.Range("B25").Value = 1.019005539842709
$ws.Range("D25").Value = 0.1960155098937406
$ws.Range("E25").Value = 0.2520926470316862
$ws.Range("F25").Value = 2.030815736048609
$ws.Range("G25").Value = 0.002506616415351283
$ws.Range("J25").Value = 0.3651879633957127
$ws.Range("K25").Value = 0.3957548502674229
$ws.Range("L25").Value = 0.1241869748498061
$ws.Range("M25").Value = 0.2078640766264321
$ws.Range("O25").Value = 5.080431773593546
